# Updates the cryptos list data (price / 1h volume change) on the active
# sheet, mirroring the daily GitHub Actions refresh of cryptos.xlsx.
#
# Column indices: A=1 (rank), B=2 (Coin), C=3 (Link), D=4 (Price), E=5 (Volume 1h)
#
# All of the "Price"/"Coin"/"Link" values in this sheet are stored as plain
# text (not numbers) even when they look numeric (e.g. "211.31", "1.00").
# Assigning a numeric-looking string straight to .Value lets Excel's COM
# layer auto-coerce it into a real number, which would change the cell's
# stored type. To keep these as text we mark the cell as Text ("@") before
# the assignment, then restore the "Normal" style afterwards so we don't
# leave a stray number-format behind while keeping the text type.

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
Set-TextValue $ws.Cells.Item(2, 4) "26.701.38"
Set-TextValue $ws.Cells.Item(2, 5) "  +0.41%  "

# Row 3: Ethereum
Set-TextValue $ws.Cells.Item(3, 4) "1.599.80"
Set-TextValue $ws.Cells.Item(3, 5) "  +0.38%  "

# Row 4: TetherUSD
Set-TextValue $ws.Cells.Item(4, 5) "  +0.23%  "

# Row 5: BNB
Set-TextValue $ws.Cells.Item(5, 4) "211.31"
Set-TextValue $ws.Cells.Item(5, 5) "  +0.29%  "

# Row 6: XRP
Set-TextValue $ws.Cells.Item(6, 5) "  -0.38%  "

# Row 7: USDC
Set-TextValue $ws.Cells.Item(7, 5) "  +0.19%  "

# Row 8: Dogecoin
Set-TextValue $ws.Cells.Item(8, 5) "  +0.49%  "

# Row 9: Cardano
Set-TextValue $ws.Cells.Item(9, 5) "  +1.09%  "

# Row 10: Solana
Set-TextValue $ws.Cells.Item(10, 4) "19.54"
Set-TextValue $ws.Cells.Item(10, 5) "  +0.99%  "

# Row 11: TRON
Set-TextValue $ws.Cells.Item(11, 4) "0.0842"
Set-TextValue $ws.Cells.Item(11, 5) "  +0.86%  "

# Row 12: WrappedliquidstakedEther2.0
Set-TextValue $ws.Cells.Item(12, 4) "1.824.23"
Set-TextValue $ws.Cells.Item(12, 5) "  +0.37%  "

# Row 13: WrappedEther
Set-TextValue $ws.Cells.Item(13, 4) "1.588.96"
Set-TextValue $ws.Cells.Item(13, 5) "  -0.17%  "

# Row 14: Polkadot
Set-TextValue $ws.Cells.Item(14, 5) "  +0.63%  "

# Row 15: Polygon
Set-TextValue $ws.Cells.Item(15, 4) "0.523"
Set-TextValue $ws.Cells.Item(15, 5) "  +0.80%  "

# Row 16: Litecoin
Set-TextValue $ws.Cells.Item(16, 5) "  +1.50%  "

# Row 17: WrappedBTC
Set-TextValue $ws.Cells.Item(17, 4) "26.679.96"

# Row 18: ShibaInu
Set-TextValue $ws.Cells.Item(18, 5) "  +3.89%  "

# Rows 19/20 swapped places (BitcoinCash <-> Dai) with updated values
Set-TextValue $ws.Cells.Item(19, 2) "Dai"
Set-TextValue $ws.Cells.Item(19, 3) "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Cells.Item(19, 4) "1.00"
Set-TextValue $ws.Cells.Item(19, 5) "  +0.21%  "

Set-TextValue $ws.Cells.Item(20, 2) "BitcoinCash"
Set-TextValue $ws.Cells.Item(20, 3) "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Cells.Item(20, 4) "209.43"
Set-TextValue $ws.Cells.Item(20, 5) "  +0.90%  "

# Row 21: Chainlink
Set-TextValue $ws.Cells.Item(21, 5) "  +3.63%  "

# Row 22: Uniswap
Set-TextValue $ws.Cells.Item(22, 5) "  +1.02%  "

# Row 23: Toncoin
Set-TextValue $ws.Cells.Item(23, 4) "2.30"
Set-TextValue $ws.Cells.Item(23, 5) "  +0.08%  "

# Row 24: Avalanche
Set-TextValue $ws.Cells.Item(24, 5) "  +1.06%  "

# Row 25: Monero
Set-TextValue $ws.Cells.Item(25, 4) "143.21"
Set-TextValue $ws.Cells.Item(25, 5) "  -1.23%  "

# Row 26: BinanceUSD
Set-TextValue $ws.Cells.Item(26, 5) "  +0.18%  "

# Row 27: Cosmos
Set-TextValue $ws.Cells.Item(27, 4) "7.11"
Set-TextValue $ws.Cells.Item(27, 5) "  +0.46%  "

# Row 28: Stellar
Set-TextValue $ws.Cells.Item(28, 5) "  +0.34%  "

# Row 30: Hedera
Set-TextValue $ws.Cells.Item(30, 5) "  +2.70%  "

# Row 31: PancakeSwap
Set-TextValue $ws.Cells.Item(31, 5) "  +0.22%  "

# Row 32: Filecoin
Set-TextValue $ws.Cells.Item(32, 5) "  +0.94%  "

# Row 33: InternetComputer(DFINITY)
Set-TextValue $ws.Cells.Item(33, 5) "  +1.87%  "

# Row 34: Maker
Set-TextValue $ws.Cells.Item(34, 4) "1.289.71"
Set-TextValue $ws.Cells.Item(34, 5) "  +0.99%  "

# Row 35: ImmutableX
Set-TextValue $ws.Cells.Item(35, 4) "0.620"
Set-TextValue $ws.Cells.Item(35, 5) "  -5.01%  "

# Row 36: HuobiToken
Set-TextValue $ws.Cells.Item(36, 5) "  +0.97%  "

# Row 37: LidoDAOToken
Set-TextValue $ws.Cells.Item(37, 5) "  +0.66%  "

# Row 38: VeChain
Set-TextValue $ws.Cells.Item(38, 5) "  +0.13%  "

# Row 39: WEMIXToken
Set-TextValue $ws.Cells.Item(39, 5) "  +15.97%  "

# Row 40: ARBITRUM
Set-TextValue $ws.Cells.Item(40, 5) "  -1.79%  "

# Row 41: FraxShare
Set-TextValue $ws.Cells.Item(41, 5) "  -0.72%  "

# Rows 42/43 swapped places (TrustWalletToken <-> MXToken) with updated values
Set-TextValue $ws.Cells.Item(42, 2) "MXToken"
Set-TextValue $ws.Cells.Item(42, 3) "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Cells.Item(42, 4) "2.19"
Set-TextValue $ws.Cells.Item(42, 5) "  -0.28%  "

Set-TextValue $ws.Cells.Item(43, 2) "TrustWalletToken"
Set-TextValue $ws.Cells.Item(43, 3) "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Cells.Item(43, 4) "0.784"
Set-TextValue $ws.Cells.Item(43, 5) "  +0.07%  "

# Row 44: Aave
Set-TextValue $ws.Cells.Item(44, 4) "63.08"
Set-TextValue $ws.Cells.Item(44, 5) "  -1.15%  "

# Row 45: RocketPoolETH
Set-TextValue $ws.Cells.Item(45, 4) "1.736.28"
Set-TextValue $ws.Cells.Item(45, 5) "  +0.49%  "

# Row 46: Quant
Set-TextValue $ws.Cells.Item(46, 4) "90.93"
Set-TextValue $ws.Cells.Item(46, 5) "  +1.81%  "

# Row 47: RenderToken
Set-TextValue $ws.Cells.Item(47, 5) "  -0.73%  "

# Row 48: Algorand
Set-TextValue $ws.Cells.Item(48, 5) "  -1.09%  "

# Row 49: Cronos
Set-TextValue $ws.Cells.Item(49, 5) "  +0.90%  "

# Row 50: USDD
Set-TextValue $ws.Cells.Item(50, 5) "  +0.08%  "

# Row 51: EnergySwap
Set-TextValue $ws.Cells.Item(51, 4) "7.38"
Set-TextValue $ws.Cells.Item(51, 5) "  -0.45%  "
